# TestMemberManagement.xlsx update
# "Added Product Management, Leather Management, Product Category, Blog"
#
# The underlying test-data rows are rotated forward to new user/record
# numbers (e.g. UserTest-55..60 -> UserTest-79..84) across the four sheets,
# the D-column numeric id sequences are bumped, a few helper columns get an
# explicit width, and the active sheet/selection bookmarks move around.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "createUser" (1)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("createUser")

$ws1.Range("A2").Value = "UserTest-79"
$ws1.Range("C2").Value = "user.testAuto079"
$ws1.Range("D2").Value = 9800000039

$ws1.Range("A3").Value = "UserTest-80"
$ws1.Range("C3").Value = "user.testAuto080"
$ws1.Range("D3").Value = 9800000040

$ws1.Range("A4").Value = "UserTest-81"
$ws1.Range("C4").Value = "user.testAuto081"
$ws1.Range("D4").Value = 9800000041

$ws1.Range("A5").Value = "UserTest-82"
$ws1.Range("C5").Value = "user.testAuto082"
$ws1.Range("D5").Value = 9800000042

$ws1.Range("A6").Value = "UserTest-83"
$ws1.Range("C6").Value = "user.testAuto083"
$ws1.Range("D6").Value = 9800000043

$ws1.Range("A7").Value = "UserTest-84"
$ws1.Range("C7").Value = "user.testAuto084"
$ws1.Range("D7").Value = 9800000044

$ws1.Columns.Item(1).ColumnWidth = 10.166666666666666

# ---------------------------------------------------------------------
# Sheet "editUser" (2)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("editUser")

$ws2.Range("A2").Value = "TestUser-66"
$ws2.Range("C2").Value = "userEdit.auto66"
$ws2.Range("D2").Value = 9800000021

$ws2.Range("A3").Value = "TestUser-67"
$ws2.Range("C3").Value = "userEdit.auto67"
$ws2.Range("D3").Value = 9800000022

$ws2.Range("A4").Value = "TestUser-68"
$ws2.Range("C4").Value = "userEdit.auto68"
$ws2.Range("D4").Value = 9800000023

$ws2.Range("A5").Value = "TestUser-69"
$ws2.Range("C5").Value = "userEdit.auto69"
$ws2.Range("D5").Value = 9800000024

$ws2.Range("A6").Value = "TestUser-70"
$ws2.Range("C6").Value = "userEdit.auto70"
$ws2.Range("D6").Value = 9800000025

$ws2.Columns.Item(1).ColumnWidth = 11.276041666666666
$ws2.Columns.Item(3).ColumnWidth = 12.053385416666666

# ---------------------------------------------------------------------
# Sheet "resetPassword" (3)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("resetPassword")

$ws3.Range("A2").Value = "TestUser-65"
$ws3.Range("C2").Value = "userPass.auto65"
$ws3.Range("D2").Value = 9800000020

$ws3.Range("A3").Value = "TestUser-66"
$ws3.Range("C3").Value = "userPass.auto66"
$ws3.Range("D3").Value = 9800000021

$ws3.Range("A4").Value = "TestUser-67"
$ws3.Range("C4").Value = "userPass.auto67"
$ws3.Range("D4").Value = 9800000022

$ws3.Range("A5").Value = "TestUser-68"
$ws3.Range("C5").Value = "userPass.auto68"
$ws3.Range("D5").Value = 9800000023

$ws3.Range("A6").Value = "TestUser-69"
$ws3.Range("C6").Value = "userPass.auto69"
$ws3.Range("D6").Value = 9800000024

$ws3.Columns.Item(1).ColumnWidth = 10.385416666666666

# ---------------------------------------------------------------------
# Sheet "checkLogin" (4)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("checkLogin")

$ws4.Range("A2").Value = "UserTest-67"
$ws4.Range("C2").Value = "user_logintest.auto67"
$ws4.Range("D2").Value = 9800000027

$ws4.Range("A3").Value = "UserTest-68"
$ws4.Range("C3").Value = "user_logintest.auto68"
$ws4.Range("D3").Value = 9800000028

$ws4.Range("A4").Value = "UserTest-69"
$ws4.Range("C4").Value = "user_logintest.auto69"
$ws4.Range("D4").Value = 9800000029

$ws4.Range("A5").Value = "UserTest-70"
$ws4.Range("C5").Value = "user_logintest.auto70"
$ws4.Range("D5").Value = 9800000030

$ws4.Range("A6").Value = "UserTest-71"
$ws4.Range("C6").Value = "user_logintest.auto71"
$ws4.Range("D6").Value = 9800000031

$ws4.Range("A7").Value = "UserTest-72"
$ws4.Range("C7").Value = "user_logintest.auto72"
$ws4.Range("D7").Value = 9800000032

# ---------------------------------------------------------------------
# Selections / active sheet.
# Order matters: the last sheet selected below becomes the active tab.
# Target state: createUser ends up active & tabSelected, checkLogin loses
# tabSelected, editUser/resetPassword selections move from C to D.
# ---------------------------------------------------------------------
$ws2.Range("D2:D6").Select() | Out-Null
$ws3.Range("D2:D6").Select() | Out-Null
$ws4.Range("D15").Select() | Out-Null
$ws1.Range("D12").Select() | Out-Null
